$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title paragraph: " - " + bold "2023"  ->  " – 2023" (not bold),
#    then add a new paragraph right after it with the disclaimer text.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)

$dashRange = $titlePara.Range.Duplicate()
$dashRange.Find.Execute(" - 2023", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dashRange.Text = " – 2023"
$dashRange.Font.Bold = 0

$titlePara.Range.InsertParagraphAfter()
$noticePara = $d.Paragraphs(2)
$noticePara.Range.Text = "(Esses softwares precisam ser instalados por você o quanto antes, no seu computador pessoal)"

# ---------------------------------------------------------------------
# 2) "PENCIL  PROTOTYPER" item: drop the grammar proofing markers that
#    wrapped the run (the text itself is unchanged).
# ---------------------------------------------------------------------
$pencilRange = $d.Content.Duplicate()
$pencilRange.Find.Execute("PENCIL  PROTOTYPER", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pencilRange.Text = "__tmp__"
$pencilRange2 = $d.Content.Duplicate()
$pencilRange2.Find.Execute("__tmp__", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pencilRange2.Text = "PENCIL  PROTOTYPER"

# ---------------------------------------------------------------------
# 3) Figma hyperlink: collapse the many spell-checked runs into a single
#    run reading "Figma: the collaborative interface design tool."
# ---------------------------------------------------------------------
$figmaRange = $d.Content.Duplicate()
$figmaRange.Find.Execute("Figma: the collaborative interface design tool.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$figmaRange.Text = "__tmp2__"
$figmaRange2 = $d.Content.Duplicate()
$figmaRange2.Find.Execute("__tmp2__", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$figmaRange2.Text = "Figma: the collaborative interface design tool."
